$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: mark part of the sentence in red.
#   "属于不同子网，跨子网通信需要默认网关的转发。而要和默认网关通信，就需要获得其MAC地址。"
#   becomes three segments with the middle clause colored red:
#     "属于不同子网，"                         (unchanged)
#     "跨子网通信需要默认网关的转发"             (red)
#     "。"                                    (unchanged)
#     "而要和默认网关通信，就需要获得其"         (red)
#     "MAC"                                   (red)
#     "地址"                                   (red)
#     "。"                                    (unchanged)
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("属于不同子网，跨子网通信需要默认网关的转发。而要和默认网关通信，就需要获得其MAC地址。", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $base = $rng.Start
    # "属于不同子网，" is 7 characters long -> red run starts right after it.
    $redStart = $base + 7
    # "。" sits right before "而要和默认网关通信" (character offset 21 within the match).
    $periodStart = $base + 21
    $periodEnd = $periodStart + 1
    # the final "。" is the very last character of the match.
    $tailPeriodStart = $rng.End - 1

    $seg1 = $d.Range($redStart, $periodStart)
    $seg1.Font.Color = 255

    $seg2 = $d.Range($periodEnd, $tailPeriodStart)
    $seg2.Font.Color = 255
}

# ---------------------------------------------------------------------
# Change 2: split "Destination: vmware_51:f1" and drop the _GoBack
# bookmark there (Word automatically relocates the single _GoBack
# bookmark, which also removes it from its old location further down).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Destination: vmware_51:f1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $splitPos = $rng2.Start + "Destination: ".Length
    $d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
}

Write-Output "stage1-2 done"
